$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates (CORREÇÕES A SEREM REALIZADAS / ARTEFATO COM ERRO columns) ---
# Row 5 (User.java comments finding): reword the "artefato com erro" text and replace the
# suggested fix with a JavaDoc recommendation.
$ws.Range("D5").Value = "O código não possui nenhum comentário descritivo útil sobre seu funcionamento"
$ws.Range("E5").Value = "Adicionar JavaDoc explicando parâmetros, retornos e exceções"

# Row 9 (connections not closed finding): replace the "Inserir close()" suggestion with a
# recommendation to use try-with-resources.
$ws.Range("E9").Value = "Usar try-with-resources para fechar ``Connection``, ``PreparedStatement`` e ``ResultSet`` automaticamente"

# The longer text in E9 now needs more vertical room.
$ws.Rows.Item(9).RowHeight = 64.5

# --- View state ---
# Scroll down a bit and move the selection to E10 (just below the table).
$ws.Range("E10").Select()
